$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 141
$ws1.Range("F4").Value = 2102
$ws1.Range("G4").Value = 60
$ws1.Range("F6").Value = 660
$ws1.Range("F7").Value = 108
$ws1.Range("F8").Value = 2088
$ws1.Range("F9").Value = 10839
$ws1.Range("F15").Value = 9057
$ws1.Range("F18").Value = 5319
$ws1.Range("F20").Value = 3369

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 141
$ws4.Range("F4").Value = 2102
$ws4.Range("G4").Value = 60
$ws4.Range("F6").Value = 660
$ws4.Range("F8").Value = 108
$ws4.Range("F9").Value = 2088
$ws4.Range("F12").Value = 10839
$ws4.Range("F18").Value = 9057
$ws4.Range("F21").Value = 5319
$ws4.Range("F23").Value = 3369
